$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4200.9165
$ws.Range("I40").Value = 3163.4443
$ws.Range("J40").Value = 7313.3335
$ws.Range("K40").Value = 3163.4443
$ws.Range("L40").Value = 7313.3335
$ws.Range("M40").Value = -2988.4443
$ws.Range("N40").Value = -7663.3335
$ws.Range("H116").Value = 3583.3333
$ws.Range("I116").Value = 3125.5
$ws.Range("K116").Value = 3125.5
$ws.Range("M116").Value = 316.5
$ws.Range("H136").Value = 148333
$ws.Range("J136").Value = 148333
$ws.Range("L136").Value = 148333
$ws.Range("N136").Value = -158533
$ws.Range("H137").Value = 1897.25
$ws.Range("I137").Value = 1299.75
$ws.Range("K137").Value = 3899.25
$ws.Range("M137").Value = -1349.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9006.77
$ws.Range("I32").Value = 7840.6665
$ws.Range("K32").Value = 7840.6665
$ws.Range("M32").Value = -7553.6665
$ws.Range("H45").Value = 2429.5833
$ws.Range("I45").Value = 1880.5
$ws.Range("J45").Value = 2978.6667
$ws.Range("K45").Value = 1880.5
$ws.Range("L45").Value = 2978.6667
$ws.Range("M45").Value = -1503.5
$ws.Range("N45").Value = -3732.6667
$ws.Range("H74").Value = 6274.8335
$ws.Range("I74").Value = 6245.0557
$ws.Range("K74").Value = 6245.0557
$ws.Range("M74").Value = -5371.0557
$ws.Range("H77").Value = 6274.8335
$ws.Range("I77").Value = 6245.0557
$ws.Range("K77").Value = 31225.2785
$ws.Range("M77").Value = -26857.2785
$ws.Range("H96").Value = 1695788.8
$ws.Range("J96").Value = 1695788.8
$ws.Range("L96").Value = 1695788.8
$ws.Range("N96").Value = -1701280.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 62646.418
$ws.Range("J74").Value = 62646.418
$ws.Range("L74").Value = 62646.418
$ws.Range("N74").Value = -64394.418
$ws.Range("H77").Value = 62646.418
$ws.Range("J77").Value = 62646.418
$ws.Range("L77").Value = 187939.254
$ws.Range("N77").Value = -196675.254
$ws.Range("H132").Value = 3707.6155
$ws.Range("I132").Value = 3097.8
$ws.Range("J132").Value = 4088.75
$ws.Range("K132").Value = 9293.400000000001
$ws.Range("L132").Value = 12266.25
$ws.Range("M132").Value = -6763.400000000001
$ws.Range("N132").Value = -17326.25
$ws.Range("H141").Value = 46873.89
$ws.Range("J141").Value = 46873.89
$ws.Range("L141").Value = 46873.89
$ws.Range("N141").Value = -57233.89

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1269.4286
$ws.Range("I132").Value = 999
$ws.Range("J132").Value = 1630
$ws.Range("K132").Value = 8991
$ws.Range("L132").Value = 14670
$ws.Range("M132").Value = -6461
$ws.Range("N132").Value = -19730

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 107.888885
$ws.Range("I9").Value = 120.2
$ws.Range("K9").Value = 120.2
$ws.Range("M9").Value = 49.8
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H47").Value = 16974.75
$ws.Range("J47").Value = 16974.75
$ws.Range("L47").Value = 16974.75
$ws.Range("N47").Value = -18110.75
$ws.Range("H63").Value = 45000
$ws.Range("J63").Value = 45000
$ws.Range("L63").Value = 45000
$ws.Range("N63").Value = -46372
$ws.Range("H66").Value = 45000
$ws.Range("J66").Value = 45000
$ws.Range("L66").Value = 135000
$ws.Range("N66").Value = -141864
$ws.Range("H70").Value = 2115.6667
$ws.Range("I70").Value = 2244.4546
$ws.Range("J70").Value = 699
$ws.Range("K70").Value = 2244.4546
$ws.Range("L70").Value = 699
$ws.Range("M70").Value = -1974.4546
$ws.Range("N70").Value = -1239
$ws.Range("H73").Value = 2115.6667
$ws.Range("I73").Value = 2244.4546
$ws.Range("J73").Value = 699
$ws.Range("K73").Value = 2244.4546
$ws.Range("L73").Value = 699
$ws.Range("M73").Value = -1308.4546
$ws.Range("N73").Value = -2571
$ws.Range("H92").Value = 6477.077
$ws.Range("J92").Value = 6996
$ws.Range("L92").Value = 6996
$ws.Range("N92").Value = -10740

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 910.5
$ws.Range("I16").Value = 895
$ws.Range("J16").Value = 988
$ws.Range("K16").Value = 895
$ws.Range("L16").Value = 988
$ws.Range("M16").Value = -725
$ws.Range("N16").Value = -1328
$ws.Range("H22").Value = 956.4286
$ws.Range("I22").Value = 965.8333
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 965.8333
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -670.8333
$ws.Range("N22").Value = -1490
$ws.Range("H27").Value = 956.4286
$ws.Range("I27").Value = 965.8333
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 965.8333
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -858.8333
$ws.Range("N27").Value = -1114
$ws.Range("H61").Value = 3809.2
$ws.Range("I61").Value = 1904.4445
$ws.Range("K61").Value = 1904.4445
$ws.Range("M61").Value = -1702.4445
$ws.Range("H68").Value = 7666.5557
$ws.Range("I68").Value = 5750
$ws.Range("K68").Value = 5750
$ws.Range("M68").Value = -5001
$ws.Range("H71").Value = 7666.5557
$ws.Range("I71").Value = 5750
$ws.Range("K71").Value = 28750
$ws.Range("M71").Value = -25006
$ws.Range("H113").Value = 3809.2
$ws.Range("I113").Value = 1904.4445
$ws.Range("K113").Value = 1904.4445
$ws.Range("M113").Value = 265.5554999999999
$ws.Range("H136").Value = 1900
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -8100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496
$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716
$ws.Range("H96").Value = 947.8570999999999
$ws.Range("I96").Value = 839.3125
$ws.Range("J96").Value = 1295.2
$ws.Range("K96").Value = 839.3125
$ws.Range("L96").Value = 1295.2
$ws.Range("M96").Value = 533.6875
$ws.Range("N96").Value = -4041.2
$ws.Range("H132").Value = 1515
$ws.Range("I132").Value = 1135
$ws.Range("K132").Value = 3405
$ws.Range("M132").Value = -875
$ws.Range("H136").Value = 3498.25
$ws.Range("I136").Value = 1999.3334
$ws.Range("K136").Value = 5998.0002
$ws.Range("M136").Value = -3448.0002
